$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while guaranteeing Excel
# keeps it as a text string rather than auto-coercing it to a
# number (which happens for values such as "1.007" or "255.61").
# We temporarily force a text number-format, assign the value, then
# restore the cells original style so no visible formatting changes.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "31.311.35"
Set-TextValue $ws.Range("E2") "  +1.65%  "
Set-TextValue $ws.Range("D3") "1.999.02"
Set-TextValue $ws.Range("E3") "  +2.26%  "
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("E4") "  +0.75%  "
Set-TextValue $ws.Range("D5") "255.61"
Set-TextValue $ws.Range("E5") "  +1.62%  "
Set-TextValue $ws.Range("D6") "0.7517"
Set-TextValue $ws.Range("E6") "  +25.55%  "
Set-TextValue $ws.Range("D7") "1.004"
Set-TextValue $ws.Range("E7") "  +0.51%  "
Set-TextValue $ws.Range("D8") "0.3429"
Set-TextValue $ws.Range("E8") "  +9.06%  "
Set-TextValue $ws.Range("D9") "27.55"
Set-TextValue $ws.Range("E9") "  +12.43%  "
Set-TextValue $ws.Range("D10") "0.07155"
Set-TextValue $ws.Range("E10") "  +4.15%  "
Set-TextValue $ws.Range("D11") "0.8345"
Set-TextValue $ws.Range("E11") "  +3.53%  "
Set-TextValue $ws.Range("D12") "0.08193"
Set-TextValue $ws.Range("E12") "  +2.85%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "2.005.94"
Set-TextValue $ws.Range("E13") "  +2.79%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D14") "100.63"
Set-TextValue $ws.Range("E14") "  -0.71%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "5.613"
Set-TextValue $ws.Range("E15") "  +4.92%  "
Set-TextValue $ws.Range("D16") "15.58"
Set-TextValue $ws.Range("E16") "  +13.33%  "
Set-TextValue $ws.Range("D17") "268.67"
Set-TextValue $ws.Range("E17") "  -4.44%  "
Set-TextValue $ws.Range("D18") "31.336.98"
Set-TextValue $ws.Range("E18") "  +1.77%  "
Set-TextValue $ws.Range("D19") "6.104"
Set-TextValue $ws.Range("E19") "  +9.06%  "
Set-TextValue $ws.Range("D20") "0.000008309"
Set-TextValue $ws.Range("E20") "  +7.97%  "
Set-TextValue $ws.Range("D21") "2.266.41"
Set-TextValue $ws.Range("E21") "  +4.20%  "
Set-TextValue $ws.Range("D22") "1.005"
Set-TextValue $ws.Range("E22") "  +0.57%  "
Set-TextValue $ws.Range("D23") "1.008"
Set-TextValue $ws.Range("E23") "  +0.89%  "
Set-TextValue $ws.Range("D24") "7.101"
Set-TextValue $ws.Range("E24") "  +6.94%  "
Set-TextValue $ws.Range("D25") "10.05"
Set-TextValue $ws.Range("E25") "  +6.16%  "
Set-TextValue $ws.Range("D26") "164.03"
Set-TextValue $ws.Range("E26") "  -0.74%  "
Set-TextValue $ws.Range("D27") "19.95"
Set-TextValue $ws.Range("E27") "  +2.04%  "
Set-TextValue $ws.Range("D28") "2.385"
Set-TextValue $ws.Range("E28") "  +14.78%  "
Set-TextValue $ws.Range("D29") "0.1345"
Set-TextValue $ws.Range("E29") "  +20.96%  "
Set-TextValue $ws.Range("D30") "1.603"
Set-TextValue $ws.Range("E30") "  +3.85%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D31") "1.384"
Set-TextValue $ws.Range("E31") "  +2.31%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.671"
Set-TextValue $ws.Range("E32") "  +4.39%  "
Set-TextValue $ws.Range("D33") "4.465"
Set-TextValue $ws.Range("E33") "  +2.75%  "
Set-TextValue $ws.Range("D34") "0.05375"
Set-TextValue $ws.Range("E34") "  +7.77%  "
Set-TextValue $ws.Range("D35") "1.296"
Set-TextValue $ws.Range("E35") "  +10.22%  "
Set-TextValue $ws.Range("D36") "0.7899"
Set-TextValue $ws.Range("E36") "  +9.39%  "
Set-TextValue $ws.Range("D37") "2.793"
Set-TextValue $ws.Range("E37") "  +2.89%  "
Set-TextValue $ws.Range("D38") "1.005"
Set-TextValue $ws.Range("E38") "  +0.56%  "
Set-TextValue $ws.Range("D39") "0.02015"
Set-TextValue $ws.Range("E39") "  +2.62%  "
Set-TextValue $ws.Range("D40") "2.912"
Set-TextValue $ws.Range("E40") "  -0.41%  "
Set-TextValue $ws.Range("D41") "85.53"
Set-TextValue $ws.Range("E41") "  +10.13%  "
Set-TextValue $ws.Range("D42") "6.835"
Set-TextValue $ws.Range("E42") "  +6.55%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D43") "0.4661"
Set-TextValue $ws.Range("E43") "  +3.17%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D44") "2.130"
Set-TextValue $ws.Range("E44") "  +5.45%  "
Set-TextValue $ws.Range("D45") "0.8579"
Set-TextValue $ws.Range("E45") "  +1.62%  "
Set-TextValue $ws.Range("D46") "105.22"
Set-TextValue $ws.Range("E46") "  +2.41%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D47") "1.004"
Set-TextValue $ws.Range("E47") "  +0.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "10.16"
Set-TextValue $ws.Range("E48") "  +1.22%  "
Set-TextValue $ws.Range("D49") "7.777"
Set-TextValue $ws.Range("E49") "  +6.81%  "
Set-TextValue $ws.Range("D50") "37.59"
Set-TextValue $ws.Range("E50") "  +4.49%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.587"
Set-TextValue $ws.Range("E51") "  +12.42%  "
